# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# immediately before the old column N ("Late"), pushing the existing
# N/O/P ("Late" / heading / "Outstanding") columns one slot to the
# right (-> O/P/Q). The sheet becomes the active tab/selection moves
# to H15, and the "NewLoanInput" sheet (previously active) is no
# longer the selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Match the width Excel gives a freshly inserted column (it inherits
# the width of the column immediately to its left - "In Advance", M).
$inheritedWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column before N; existing N/O/P shift to O/P/Q.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $inheritedWidth

# Make "Repayment schedule" the active sheet/tab, with H15 selected.
$ws.Activate()
$ws.Range("H15").Select()
